# Update cryptocurrency price/volume data to the latest GitHub Actions scrape.
# Cells in column D that now look like plain numbers are forced back to
# Text (NumberFormat "@") before the assignment so they are stored the same
# way as the other already-numeric-looking "Price" strings in this column
# (e.g. "1.00", "413.35") instead of being auto-converted to real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.066.24'
$ws.Range('E2').Value = '  +7.26%  '
# Row 3
$ws.Range('D3').Value = '3.356.26'
$ws.Range('E3').Value = '  +3.58%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.16%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '413.75'
$ws.Range('E5').Value = '  +4.80%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.65'
$ws.Range('E6').Value = '  +4.99%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.587'
$ws.Range('E7').Value = '  +3.71%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.08%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.640'
$ws.Range('E9').Value = '  +3.53%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.85'
$ws.Range('E10').Value = '  +2.41%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0995'
$ws.Range('E11').Value = '  +2.55%  '
# Row 12
$ws.Range('E12').Value = '  +1.24%  '
# Row 13
$ws.Range('D13').Value = '3.887.39'
$ws.Range('E13').Value = '  +3.36%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.43'
$ws.Range('E14').Value = '  +3.60%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.88'
$ws.Range('E15').Value = '  +4.62%  '
# Row 16
$ws.Range('D16').Value = '3.343.72'
$ws.Range('E16').Value = '  +3.27%  '
# Row 17
$ws.Range('E17').Value = '  +1.90%  '
# Row 18
$ws.Range('D18').Value = '60.747.46'
$ws.Range('E18').Value = '  +6.92%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.78'
$ws.Range('E19').Value = '  +1.42%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.39'
$ws.Range('E20').Value = '  +2.43%  '
# Row 21
$ws.Range('E21').Value = '  +5.64%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.01'
$ws.Range('E22').Value = '  +0.02%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '303.43'
$ws.Range('E23').Value = '  +1.50%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.07'
$ws.Range('E24').Value = '  +1.62%  '
# Row 25
$ws.Range('E25').Value = '  +2.19%  '
# Row 26
$ws.Range('E26').Value = '  +3.20%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.48'
$ws.Range('E27').Value = '  +2.07%  '
# Row 28
$ws.Range('E28').Value = '  +6.78%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.97'
$ws.Range('E29').Value = '  +1.32%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.66'
$ws.Range('E30').Value = '  +5.80%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.64'
$ws.Range('E31').Value = '  +24.60%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.114'
$ws.Range('E32').Value = '  +4.70%  '
# Row 33
$ws.Range('E33').Value = '  +4.29%  '
# Row 34
$ws.Range('E34').Value = '  -0.02%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '39.12'
$ws.Range('E35').Value = '  +3.58%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0510'
$ws.Range('E36').Value = '  +5.68%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.51'
$ws.Range('E37').Value = '  +1.61%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.13'
$ws.Range('E38').Value = '  +2.39%  '
# Row 39
$ws.Range('E39').Value = '  +0.15%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.42'
$ws.Range('E40').Value = '  -1.69%  '
# Row 41
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '136.84'
$ws.Range('E41').Value = '  +2.06%  '
# Row 42
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.297'
$ws.Range('E42').Value = '  +3.68%  '
# Row 43
$ws.Range('E43').Value = '  +2.93%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.92'
$ws.Range('E44').Value = '  +0.69%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.97'
$ws.Range('E45').Value = '  +0.06%  '
# Row 46
$ws.Range('E46').Value = '  -1.08%  '
# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.68'
$ws.Range('E47').Value = '  +3.45%  '
# Row 48
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.25'
$ws.Range('E48').Value = '  +8.82%  '
# Row 49
$ws.Range('D49').Value = '2.174.57'
$ws.Range('E49').Value = '  +1.67%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.41'
$ws.Range('E50').Value = '  +1.05%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.98'
$ws.Range('E51').Value = '  -1.55%  '
